$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new data row at row 15 (shifts existing rows 15-45 down to 16-46),
# mirroring the weekly update described in the commit message ("Fruta / hortaliza, semanal").
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 4
$ws.Cells.Item(15, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(15, 3).Value = "Los Lagos"
$ws.Cells.Item(15, 4).Value = 44665
$ws.Cells.Item(15, 5).Value = 10
$ws.Cells.Item(15, 6).Value = 100112030
$ws.Cells.Item(15, 7).Value = "Poroto granado"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 35
$ws.Cells.Item(15, 11).Value = 28000
$ws.Cells.Item(15, 12).Value = 28000
$ws.Cells.Item(15, 13).Value = 28000
$ws.Cells.Item(15, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Región Metropolitana"
$ws.Cells.Item(15, 16).Value = 1120
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
